$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 289-290, pushing the existing rows 289.. down to 291..
$ws.Range("A289:A290").EntireRow.Insert()

# Fill in the new row 289 (Lane Late, Primera)
$ws.Range("A289").Value2 = 11
$ws.Range("B289").Value2 = "Vega Monumental Concepción"
$ws.Range("C289").Value2 = "Bíobío"
$ws.Range("D289").Value2 = 44876
$ws.Range("E289").Value2 = 8
$ws.Range("F289").Value2 = "Fruta"
$ws.Range("G289").Value2 = 100102
$ws.Range("H289").Value2 = "Cítricos"
$ws.Range("I289").Value2 = 100102005
$ws.Range("J289").Value2 = "Naranja"
$ws.Range("K289").Value2 = "Lane Late"
$ws.Range("L289").Value2 = "Primera"
$ws.Range("M289").Value2 = 500
$ws.Range("N289").Value2 = 9500
$ws.Range("O289").Value2 = 10000
$ws.Range("P289").Value2 = 9750
$ws.Range("Q289").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R289").Value2 = "Región de O'Higgins"
$ws.Range("S289").Value2 = 650
$ws.Range("T289").Value2 = 15

# Fill in the new row 290 (Valencia, Primera)
$ws.Range("A290").Value2 = 11
$ws.Range("B290").Value2 = "Vega Monumental Concepción"
$ws.Range("C290").Value2 = "Bíobío"
$ws.Range("D290").Value2 = 44876
$ws.Range("E290").Value2 = 8
$ws.Range("F290").Value2 = "Fruta"
$ws.Range("G290").Value2 = 100102
$ws.Range("H290").Value2 = "Cítricos"
$ws.Range("I290").Value2 = 100102005
$ws.Range("J290").Value2 = "Naranja"
$ws.Range("K290").Value2 = "Valencia"
$ws.Range("L290").Value2 = "Primera"
$ws.Range("M290").Value2 = 400
$ws.Range("N290").Value2 = 9000
$ws.Range("O290").Value2 = 10000
$ws.Range("P290").Value2 = 9625
$ws.Range("Q290").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R290").Value2 = "Región de O'Higgins"
$ws.Range("S290").Value2 = 642
$ws.Range("T290").Value2 = 15
